# Update: hinge calculation for elevator and rudder plus engine loads
# stored correctly inside the structure variable.
#
# - Corrects three engine-load figures (Takeoff_power, Max_Continous_power,
#   Engine_mount_mass/accessories/spinner) that had been stored using the
#   wrong (unconverted) magnitude.
# - Appends ten new rows (121-130) describing the vertical tail / rudder
#   geometry and hinge-moment coefficients needed for the rudder hinge
#   calculation, mirroring the existing elevator block above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected engine figures -------------------------------------------------
$ws.Range("B93").Value  = 11.185499999999999   # Takeoff_power [kW]
$ws.Range("B95").Value  = 9.3212480000000006   # Max_Continous_power [kW]
$ws.Range("B101").Value = 1.4                  # Engine_mount_mass [kg]
$ws.Range("B102").Value = 18.5                 # Engine_accessories [kg]
$ws.Range("B103").Value = 4.5                  # Propeller_spinner [kg]

# --- New vertical tail / rudder geometry + hinge data (rows 121-130) ---------
$ws.Range("A121").Value = "S_vertical"
$ws.Range("B121").Value = 0.1022
$ws.Range("C121").Value = "m^2"

$ws.Range("A122").Value = "chord_vertical"
$ws.Range("B122").Value = 0.3136
$ws.Range("C122").Value = "m"

$ws.Range("A123").Value = "S_rudder"
$ws.Range("B123").Value = 0.0381
$ws.Range("C123").Value = "m^2"

$ws.Range("A124").Value = "chord_rudder"
$ws.Range("B124").Value = 0.14
$ws.Range("C124").Value = "m"

$ws.Range("A125").Value = "chord_ratio_rudder_cf_c"
$ws.Range("B125").Value = 0.35
$ws.Range("C125").Value = "Non dimensional"

$ws.Range("A126").Value = "overhang_rudder"
$ws.Range("B126").Value = 0.12
$ws.Range("C126").Value = "Non dimensional"

$ws.Range("A127").Value = "span_ratio_rudder"
$ws.Range("B127").Value = 0.8
$ws.Range("C127").Value = "Non dimensional"

$ws.Range("A128").Value = "max_deflection_rudder"
$ws.Range("B128").Value = 25
$ws.Range("C128").Value = "degrees"

$ws.Range("A129").Value = "Chdeltarudder"
$ws.Range("B129").Value = -0.4538
$ws.Range("C129").Value = "1/rad"

$ws.Range("A130").Value = "Chalfarudder"
$ws.Range("B130").Value = -0.0024
$ws.Range("C130").Value = "1/rad"

# --- Restore the view: scroll/selection now centred on the new rows ---------
$ws.Range("B104").Select()
